# Apply "Warmachine games - list all or by user. None faction added" edit.
$wb = $excel.ActiveWorkbook

# Worksheets in tab order:
#  1: M-1 Tasks
#  2: M-0.5
#  3: M0 - Account Mgmt   (sheet3.xml)
#  4: M1 - Game Data      (sheet4.xml)
#  5: Links
#  6: Notes
$wsAccount = $wb.Worksheets.Item(3)
$wsGame    = $wb.Worksheets.Item(4)

# ---------------------------------------------------------------------
# M1 - Game Data sheet content updates
# ---------------------------------------------------------------------

# D10: expand note about restricting edits to the owning user
$wsGame.Range("D10").Value = "Need to restrict this to the owning user -> 1) Edit only on own controller (2) Compare ids in controller"

# A14: fix typo ("entried" -> "entries") - now reads the same as the (removed) duplicate row
$wsGame.Range("A14").Value = "User can view all game entries by ANOTHER user"

# D14: replace stale note with a new two-run rich text note (plain + bold warning)
$prefix = "Link to profile page. --- Techincally done, but need way to pass a userid. "
$boldPart = "DO USER PROFILE VIEW NEXT"
$d14 = $wsGame.Range("D14")
$d14.Value = $prefix + $boldPart
$d14.Characters($prefix.Length + 1, $boldPart.Length).Font.Bold = $true

# C13: new "Done" status cell (Good style, matches other status cells in column C)
$wsGame.Range("C13").Value = "Done"
$wsGame.Range("C13").Style = "Good"

# D19: new note about needing a data migration, bold styled like the other empty bold note cell (D18)
$wsGame.Range("D19").Value = "FIX THIS! SERIOUSLY! Will need data migration"
$wsGame.Range("D19").Font.Bold = $true

# ---------------------------------------------------------------------
# View / selection / active sheet updates
# ---------------------------------------------------------------------

# "M0 - Account Mgmt" keeps a selection on A19:D19 but is no longer the active tab
$wsAccount.Range("A19:D19").Select()

# "M1 - Game Data" becomes the active sheet/tab, with D19 selected
$wsGame.Activate()
$wsGame.Range("D19").Select()
